$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2..11 (columns A..J)
# A = serial date/time (numeric), B..G,I numeric, H/J text (may be blank)
$data = @(
    @(45846,                 2025, 32, 14.07, 91.88, 0,      0.68, "ESE", 0, ""),
    @(45846.04166666666,     2025, 32, 13.94, 91.95, 0,      0,    "",    0, ""),
    @(45846.08333333334,     2025, 32, 14.13, 91.93000000000001, 0, 3.23, "NW",  0, ""),
    @(45846.125,             2025, 32, 14.26, 90.53, 0,      6.21, "E",   0, ""),
    @(45846.16666666666,     2025, 32, 13.69, 91.88, 0,      3.99, "ESE", 0, ""),
    @(45846.20833333334,     2025, 32, 13.78, 92.37, 0,      2.19, "NW",  0, ""),
    @(45846.25,              2025, 32, 13.69, 92.51000000000001, 0.88, 3.6, "SE", 0, ""),
    @(45846.29166666666,     2025, 32, 13.64, 92.79000000000001, 24.61, 2.89, "E", 0, ""),
    @(45846.33333333334,     2025, 32, 13.8,  92.93000000000001, 62.7, 1.6, "E",  0, ""),
    @(45876.37517238504,     2025, 28, 14,    93.12, 110.74, 2.16, "NE", 0, "09:00:14")
)

# Capture the date/time number format already used on A2 so every row in
# column A keeps the same display format (matches style index "2").
$dateFormat = $ws.Cells.Item(2, 1).NumberFormat

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws.Cells.Item($rowIndex, 8).Value = $row[7]
    $ws.Cells.Item($rowIndex, 9).Value = $row[8]
    $ws.Cells.Item($rowIndex, 10).Value = $row[9]
    $rowIndex++
}
